# The document ends with two empty paragraphs right before the section
# properties. The last of those two empty paragraphs gets turned into a
# hyperlink (displaying its own URL as text), matching the style used by
# every other hyperlink already present in the document.

$d = $word.ActiveDocument

$lastParaIndex = $d.Paragraphs.Count
$targetPara = $d.Paragraphs.Item($lastParaIndex)
$targetRange = $targetPara.Range

$url = "https://www.google.com/url?sa=i&url=https%3A%2F%2Fwww.vectorstock.com%2Froyalty-free-vector%2Fretro-spaceship-pixel-art-game-rocket-at-night-vector-26751070&psig=AOvVaw2WxfEWMx23tqtM_h8ekVOH&ust=1644487208537000&source=images&cd=vfe&ved=0CAwQjhxqFwoTCMCh97Cu8vUCFQAAAAAdAAAAABAD"

$d.Hyperlinks.Add($targetRange, $url, [Type]::Missing, [Type]::Missing, $url, [Type]::Missing)
